$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 3
$ws.Range("G2").Value = 0.717901
$ws.Range("H2").Value = 2.153703
$ws.Range("I2").Value = 0.0380297505351077
$ws.Range("J2").Value = 0.0380297505351077
$ws.Range("K2").Value = 3
$ws.Range("M2").Value = 2.228108666666667
$ws.Range("N2").Value = 6.684326
$ws.Range("O2").Value = 0.5687849952918405
$ws.Range("P2").Value = 0.5687849952918405
$ws.Range("Q2").Value = 1.599561439908667
$ws.Range("R2").Value = 14.396052959178
$ws.Range("S2").Value = 0.0216307514790611
$ws.Range("T2").Value = 0.0216307514790611

$ws.Range("E3").Value = 3
$ws.Range("G3").Value = 0.717901
$ws.Range("H3").Value = 2.153703
$ws.Range("I3").Value = 0.0380297505351077
$ws.Range("J3").Value = 0.0380297505351077
$ws.Range("K3").Value = 3
$ws.Range("M3").Value = 0.1124773333333333
$ws.Range("N3").Value = 0.337432
$ws.Range("O3").Value = 0.02871288122861097
$ws.Range("P3").Value = 0.02871288122861097
$ws.Range("Q3").Value = 0.08074759007733333
$ws.Range("R3").Value = 0.726728310696
$ws.Range("S3").Value = 0.001091943710268252
$ws.Range("T3").Value = 0.001091943710268252

$ws.Range("E4").Value = 3
$ws.Range("G4").Value = 0.717901
$ws.Range("H4").Value = 2.153703
$ws.Range("I4").Value = 0.0380297505351077
$ws.Range("J4").Value = 0.0380297505351077
$ws.Range("K4").Value = 3
$ws.Range("M4").Value = 1.576726666666667
$ws.Range("N4").Value = 4.73018
$ws.Range("O4").Value = 0.4025021234795486
$ws.Range("P4").Value = 0.4025021234795487
$ws.Range("Q4").Value = 1.131933650726667
$ws.Range("R4").Value = 10.18740285654
$ws.Range("S4").Value = 0.01530705534577835
$ws.Range("T4").Value = 0.01530705534577835

$ws.Range("E5").Value = 3
$ws.Range("G5").Value = 13.91986866666667
$ws.Range("H5").Value = 41.75960600000001
$ws.Range("I5").Value = 0.7373845876726675
$ws.Range("J5").Value = 0.7373845876726675
$ws.Range("K5").Value = 3
$ws.Range("M5").Value = 2.228108666666667
$ws.Range("N5").Value = 6.684326
$ws.Range("O5").Value = 0.5687849952918405
$ws.Range("P5").Value = 0.5687849952918405
$ws.Range("Q5").Value = 31.01498001506178
$ws.Range("R5").Value = 279.134820135556
$ws.Range("S5").Value = 0.4194132892276739
$ws.Range("T5").Value = 0.4194132892276739

$ws.Range("E6").Value = 3
$ws.Range("G6").Value = 13.91986866666667
$ws.Range("H6").Value = 41.75960600000001
$ws.Range("I6").Value = 0.7373845876726675
$ws.Range("J6").Value = 0.7373845876726675
$ws.Range("K6").Value = 3
$ws.Range("M6").Value = 0.1124773333333333
$ws.Range("N6").Value = 0.337432
$ws.Range("O6").Value = 0.02871288122861097
$ws.Range("P6").Value = 0.02871288122861097
$ws.Range("Q6").Value = 1.565669707976889
$ws.Range("R6").Value = 14.091027371792
$ws.Range("S6").Value = 0.02117243608565357
$ws.Range("T6").Value = 0.02117243608565358

$ws.Range("E7").Value = 3
$ws.Range("G7").Value = 13.91986866666667
$ws.Range("H7").Value = 41.75960600000001
$ws.Range("I7").Value = 0.7373845876726675
$ws.Range("J7").Value = 0.7373845876726675
$ws.Range("K7").Value = 3
$ws.Range("M7").Value = 1.576726666666667
$ws.Range("N7").Value = 4.73018
$ws.Range("O7").Value = 0.4025021234795486
$ws.Range("P7").Value = 0.4025021234795487
$ws.Range("Q7").Value = 21.94782812323111
$ws.Range("R7").Value = 197.53045310908
$ws.Range("S7").Value = 0.2967988623593401
$ws.Range("T7").Value = 0.2967988623593401

$ws.Range("E8").Value = 3
$ws.Range("G8").Value = 4.239582666666666
$ws.Range("H8").Value = 12.718748
$ws.Range("I8").Value = 0.2245856617922248
$ws.Range("J8").Value = 0.2245856617922248
$ws.Range("K8").Value = 3
$ws.Range("M8").Value = 2.228108666666667
$ws.Range("N8").Value = 6.684326
$ws.Range("O8").Value = 0.5687849952918405
$ws.Range("P8").Value = 0.5687849952918405
$ws.Range("Q8").Value = 9.446250882649778
$ws.Range("R8").Value = 85.016257943848
$ws.Range("S8").Value = 0.1277409545851055
$ws.Range("T8").Value = 0.1277409545851055

$ws.Range("E9").Value = 3
$ws.Range("G9").Value = 4.239582666666666
$ws.Range("H9").Value = 12.718748
$ws.Range("I9").Value = 0.2245856617922248
$ws.Range("J9").Value = 0.2245856617922248
$ws.Range("K9").Value = 3
$ws.Range("M9").Value = 0.1124773333333333
$ws.Range("N9").Value = 0.337432
$ws.Range("O9").Value = 0.02871288122861097
$ws.Range("P9").Value = 0.02871288122861097
$ws.Range("Q9").Value = 0.4768569527928888
$ws.Range("R9").Value = 4.291712575136
$ws.Range("S9").Value = 0.006448501432689144
$ws.Range("T9").Value = 0.006448501432689145

$ws.Range("E10").Value = 3
$ws.Range("G10").Value = 4.239582666666666
$ws.Range("H10").Value = 12.718748
$ws.Range("I10").Value = 0.2245856617922248
$ws.Range("J10").Value = 0.2245856617922248
$ws.Range("K10").Value = 3
$ws.Range("M10").Value = 1.576726666666667
$ws.Range("N10").Value = 4.73018
$ws.Range("O10").Value = 0.4025021234795486
$ws.Range("P10").Value = 0.4025021234795487
$ws.Range("Q10").Value = 6.68466304607111
$ws.Range("R10").Value = 60.16196741464
$ws.Range("S10").Value = 0.09039620577443021
$ws.Range("T10").Value = 0.09039620577443024

